# Weekly update: insert one new daily price record for Femacal de La Calera - Mango.
# A new row is inserted at row 257 (pushing the existing rows 257:370 down to 258:371),
# and the new row is populated with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 257, shifting existing rows down.
$ws.Rows.Item(257).Insert()

# Populate the newly inserted row with the new data record.
$ws.Cells.Item(257, 1).Value  = 3
$ws.Cells.Item(257, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(257, 3).Value  = "Coquimbo"
$ws.Cells.Item(257, 4).Value  = 44636
$ws.Cells.Item(257, 5).Value  = 5
$ws.Cells.Item(257, 6).Value  = "Fruta"
$ws.Cells.Item(257, 7).Value  = 100108
$ws.Cells.Item(257, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(257, 9).Value  = 100108002
$ws.Cells.Item(257, 10).Value = "Mango"
$ws.Cells.Item(257, 11).Value = "Sin especificar"
$ws.Cells.Item(257, 12).Value = "Primera"
$ws.Cells.Item(257, 13).Value = 456
$ws.Cells.Item(257, 14).Value = 7000
$ws.Cells.Item(257, 15).Value = 7000
$ws.Cells.Item(257, 16).Value = 7000
$ws.Cells.Item(257, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(257, 18).Value = "Perú"
$ws.Cells.Item(257, 19).Value = 1750
$ws.Cells.Item(257, 20).Value = 4
